# Flujo originación con compra cartera hasta analisis
#
# Updates the sample-data row on the "OriginacionDigiCredito" sheet with a
# new test case, and nudges the saved window/selection state to match.

$wb = $excel.ActiveWorkbook

# --- Target worksheet: third tab, "OriginacionDigiCredito" ---------------
$ws = $wb.Worksheets.Item("OriginacionDigiCredito")
$ws.Activate()

# --- Update the single data row (row 2) with the new test values ---------
# (write order matters for shared-string table placement on save)
$ws.Range("AX2").Value = '"86313"'       # NumRadicadoCredito
$ws.Range("B2").Value = '"52912399"'     # Cedula
$ws.Range("P2").Value = '"RODRIGUEZ"'    # Papellido
$ws.Range("Q2").Value = '"GONZALEZ"'     # Sapellido
$ws.Range("R2").Value = '"25/Nov/2021"'  # fechaActual
$ws.Range("E2").Value = '"50"'           # Plazo
$ws.Range("G2").Value = '"50"'           # DiasHabilesIntereses
$ws.Range("I2").Value = '"250000"'       # descLey
$ws.Range("H2").Value = '"8500000"'      # Ingresos
$ws.Range("J2").Value = '"300000"'       # descNomina
$ws.Range("N2").Value = '"ALEYDA"'       # Pnombre

# --- Selection / scroll position on the active sheet ----------------------
$ws.Range("F6").Select()

# --- Window position on the workbook ---------------------------------------
$win = $excel.Windows.Item(1)
$win.Left = -120
$win.Top = -120
